$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("personal_info")

$ws.Range("A1").Value = "23r"
$ws.Range("B1").Value = "23r"
